$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header G1
$ws.Range("G1").Value = "S Tag"

# Rows where Speaker (column D) changes from "RBD" to "T"
$rbdRows = @(2,3,4,6,9,11,13,15,17,18,21,22,23,24,27,29,30,32,33,35,36,37,39,46,47,48)
foreach ($r in $rbdRows) {
    $ws.Cells.Item($r, 4).Value = "T"
}

# Rows where Speaker (column D) changes from "Student" to "S"
$studentRows = @(7,12)
foreach ($r in $studentRows) {
    $ws.Cells.Item($r, 4).Value = "S"
}

# Rows where Speaker (column D) changes from "Students" to "Ss"
$studentsRows = @(5)
foreach ($r in $studentsRows) {
    $ws.Cells.Item($r, 4).Value = "Ss"
}
